$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cl = "Dependent on Cloud Infrastructure Migration milestone completion"
$it = "Critical action for Information Technology success"

# Row 8
$ws.Range("E8").Value = "Chief Technology Officer"
$ws.Range("I8").Value = $cl
$ws.Range("J8").Value = $it

# Row 9
$ws.Range("E9").Value = "IT Managers"
$ws.Range("I9").Value = $cl
$ws.Range("J9").Value = $it

# Row 10
$ws.Range("E10").Value = "DevOps Engineers"
$ws.Range("I10").Value = $cl
$ws.Range("J10").Value = $it

# Row 11
$ws.Range("E11").Value = "System Administrators"
$ws.Range("I11").Value = $cl
$ws.Range("J11").Value = $it

# Row 12 (owner unchanged)
$ws.Range("I12").Value = $cl
$ws.Range("J12").Value = $it

# Row 13 (owner unchanged)
$ws.Range("I13").Value = $cl
$ws.Range("J13").Value = $it

# Row 14
$ws.Range("E14").Value = "Chief Technology Officer"
$ws.Range("I14").Value = $cl
$ws.Range("J14").Value = $it

# Row 15
$ws.Range("E15").Value = "IT Managers"
$ws.Range("I15").Value = $cl
$ws.Range("J15").Value = $it

# Row 16
$ws.Range("E16").Value = "DevOps Engineers"
$ws.Range("I16").Value = $cl
$ws.Range("J16").Value = $it

# Row 17
$ws.Range("E17").Value = "System Administrators"
$ws.Range("I17").Value = $cl
$ws.Range("J17").Value = $it
